$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.192.15"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.859.98"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.71"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07819"
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3108"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.90"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("D12").Value = "1.859.20"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "92.72"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.123"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6908"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.557"
$ws.Range("E16").Value = "  +2.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008433"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "29.207.39"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.83"
$ws.Range("D20").Value = "2.110.52"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.92"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.597"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1533"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.87"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.890"
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.55"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  +4.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.277"
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.249"
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05223"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7569"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.873"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.707"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "1.221.59"
$ws.Range("E39").Value = "  -4.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.722"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9019"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.09"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.813"
$ws.Range("E43").Value = "  -9.37%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "2.009.96"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.70"
$ws.Range("E47").Value = "  -11.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5179"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.516"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.035"
$ws.Range("E51").Value = "  -0.93%  "
